$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 718; existing rows 718-759 (2026/12/29 .. 2027/01/05)
# shift down to 719-760.
$ws.Rows.Item(718).Insert()

# Column A holds date-like text (e.g. "2026/12/29") stored as plain text rather than
# a real Excel date. Assigning a date-looking string straight to a General-formatted
# cell's .Value gets auto-parsed into a date serial by the COM layer, so instead write
# it as a text formula into a scratch cell far outside the used range and copy only the
# computed (text) VALUE over - this avoids flipping the destination cell's number
# format / style the way forcing NumberFormat="@" would.
$scratch = $ws.Cells.Item(2000, 1)
$scratch.Formula = '="2026/01/29"'
$scratch.Copy()
$ws.Cells.Item(718, 1).PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Cells.Item(718, 2).Value = "木"
$ws.Cells.Item(718, 3).Value = 17
$ws.Cells.Item(718, 4).Value = 20
